# Generate Report for Handoff
# Adds a new tracked file (b03038be-...) to the localization-status workbook:
# a new row in "Overview", "zh-cn" and "de-de" mirroring the existing
# dc53f774-... row, but for the newly handed-off file.

$wb = $excel.ActiveWorkbook

$newFileBare = "b03038be-c130-46b1-8a4d-626d0d4123a9ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newFileE2E  = "e2e\b03038be-c130-46b1-8a4d-626d0d4123a9ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e8ce59a0598f47a4d7c5849cb7e032c18e1b1302/e2e/$newFileBare"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Overview sheet: new row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newFileBare
$wsOverview.Range("B3").Value = $newFileE2E
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-05 00:32:35"
$wsOverview.Range("G3").NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkTarget, "", "", $newFileE2E) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet: new row 3
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = $newFileBare
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "b03038be-c130-46b1-8a4d-626d0d4123a9oooooooooooooooooooooooooooooooooooooooo.44bffa351275b969317fa7dd73699ff247030595.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-05 00:32:31"
$wsZhCn.Range("H3").NumberFormat = $dateFmt
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = $dateFmt
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $hyperlinkTarget, "", "", $newFileBare) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet: new row 3
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = $newFileBare
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "b03038be-c130-46b1-8a4d-626d0d4123a9oooooooooooooooooooooooooooooooooooooooo.44bffa351275b969317fa7dd73699ff247030595.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-05 00:32:35"
$wsDeDe.Range("H3").NumberFormat = $dateFmt
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = $dateFmt
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $hyperlinkTarget, "", "", $newFileBare) | Out-Null
